$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AD hold the data to rotate/swap; column A (row index) stays put.
$firstCol = 2   # B
$lastCol  = 30  # AD

# Rows 298, 300, 302 form a 3-cycle: 298 <- 302, 300 <- 298(old), 302 <- 300(old)
# Rows 303, 304 are swapped.

# Capture "before" snapshots for the rows we need (as arrays of values B..AD).
function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += , ($ws.Cells.Item($row, $c).Value2)
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

$row298 = Get-RowValues 298
$row300 = Get-RowValues 300
$row302 = Get-RowValues 302
$row303 = Get-RowValues 303
$row304 = Get-RowValues 304

# Apply rotation: new298 = old302, new300 = old298, new302 = old300
Set-RowValues 298 $row302
Set-RowValues 300 $row298
Set-RowValues 302 $row300

# Apply swap: new303 = old304, new304 = old303
Set-RowValues 303 $row304
Set-RowValues 304 $row303
